# Updates the cryptos price/volume snapshot to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as text, even when the string looks numeric,
# matching the source data's inlineStr cell type, then restore the
# cell's default (unstyled) appearance so only the value changes.
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Cells.Item(2, 4) "29.150.46"
$ws.Cells.Item(2, 5).Value = "  -1.01%  "

# Row 3
Set-TextCell $ws.Cells.Item(3, 4) "1.860.91"
$ws.Cells.Item(3, 5).Value = "  -0.86%  "

# Row 5
Set-TextCell $ws.Cells.Item(5, 4) "0.7080"
$ws.Cells.Item(5, 5).Value = "  -0.90%  "

# Row 6
Set-TextCell $ws.Cells.Item(6, 4) "240.79"
$ws.Cells.Item(6, 5).Value = "  -0.53%  "

# Row 7
Set-TextCell $ws.Cells.Item(7, 4) "1.000"
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.92%  "

# Row 9
Set-TextCell $ws.Cells.Item(9, 4) "0.07641"
$ws.Cells.Item(9, 5).Value = "  -2.67%  "

# Row 10
Set-TextCell $ws.Cells.Item(10, 4) "24.63"

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.96%  "

# Row 12
Set-TextCell $ws.Cells.Item(12, 4) "1.865.58"
$ws.Cells.Item(12, 5).Value = "  -1.40%  "

# Row 13
Set-TextCell $ws.Cells.Item(13, 4) "5.175"
$ws.Cells.Item(13, 5).Value = "  -1.93%  "

# Row 14
Set-TextCell $ws.Cells.Item(14, 4) "0.7088"
$ws.Cells.Item(14, 5).Value = "  -2.66%  "

# Row 15
Set-TextCell $ws.Cells.Item(15, 4) "91.09"
$ws.Cells.Item(15, 5).Value = "  +0.14%  "

# Row 16
Set-TextCell $ws.Cells.Item(16, 4) "29.162.12"
$ws.Cells.Item(16, 5).Value = "  -1.09%  "

# Row 17
Set-TextCell $ws.Cells.Item(17, 4) "5.908"
$ws.Cells.Item(17, 5).Value = "  -0.60%  "

# Row 18
Set-TextCell $ws.Cells.Item(18, 4) "242.31"
$ws.Cells.Item(18, 5).Value = "  -2.01%  "

# Row 19
Set-TextCell $ws.Cells.Item(19, 4) "0.000007809"
$ws.Cells.Item(19, 5).Value = "  -0.87%  "

# Row 20
Set-TextCell $ws.Cells.Item(20, 4) "2.115.67"
$ws.Cells.Item(20, 5).Value = "  -1.04%  "

# Row 21
Set-TextCell $ws.Cells.Item(21, 4) "13.06"
$ws.Cells.Item(21, 5).Value = "  -1.83%  "

# Row 22
Set-TextCell $ws.Cells.Item(22, 4) "0.9998"
$ws.Cells.Item(22, 5).Value = "  -0.03%  "

# Row 23
Set-TextCell $ws.Cells.Item(23, 4) "7.835"
$ws.Cells.Item(23, 5).Value = "  -1.89%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.06%  "

# Row 25
Set-TextCell $ws.Cells.Item(25, 4) "0.1588"
$ws.Cells.Item(25, 5).Value = "  -0.02%  "

# Row 26
Set-TextCell $ws.Cells.Item(26, 4) "163.17"
$ws.Cells.Item(26, 5).Value = "  -0.41%  "

# Row 27
Set-TextCell $ws.Cells.Item(27, 4) "8.927"
$ws.Cells.Item(27, 5).Value = "  -0.96%  "

# Row 28
Set-TextCell $ws.Cells.Item(28, 4) "18.44"
$ws.Cells.Item(28, 5).Value = "  +0.70%  "

# Row 29
Set-TextCell $ws.Cells.Item(29, 4) "1.496"
$ws.Cells.Item(29, 5).Value = "  +0.14%  "

# Row 30
Set-TextCell $ws.Cells.Item(30, 4) "1.314"
$ws.Cells.Item(30, 5).Value = "  -3.72%  "

# Row 31
Set-TextCell $ws.Cells.Item(31, 4) "4.395"
$ws.Cells.Item(31, 5).Value = "  +0.35%  "

# Row 32
Set-TextCell $ws.Cells.Item(32, 4) "4.201"
$ws.Cells.Item(32, 5).Value = "  +1.70%  "

# Row 33
Set-TextCell $ws.Cells.Item(33, 4) "0.05126"
$ws.Cells.Item(33, 5).Value = "  -3.52%  "

# Row 34
Set-TextCell $ws.Cells.Item(34, 4) "0.7970"
$ws.Cells.Item(34, 5).Value = "  +9.70%  "

# Row 35
Set-TextCell $ws.Cells.Item(35, 4) "1.906"
$ws.Cells.Item(35, 5).Value = "  -1.57%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -3.32%  "

# Row 37
Set-TextCell $ws.Cells.Item(37, 4) "2.680"
$ws.Cells.Item(37, 5).Value = "  +0.22%  "

# Row 38
Set-TextCell $ws.Cells.Item(38, 4) "0.01838"
$ws.Cells.Item(38, 5).Value = "  -1.64%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -1.06%  "

# Row 40
Set-TextCell $ws.Cells.Item(40, 4) "1.164.80"
$ws.Cells.Item(40, 5).Value = "  -6.19%  "

# Row 41
Set-TextCell $ws.Cells.Item(41, 4) "6.170"
$ws.Cells.Item(41, 5).Value = "  +0.23%  "

# Row 42 - was TrustWalletToken, now Aave (name/link swap with new price data)
$ws.Cells.Item(42, 2).Value = "Aave"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Cells.Item(42, 4) "72.82"
$ws.Cells.Item(42, 5).Value = "  -2.01%  "

# Row 43 - was Aave, now TrustWalletToken (name/link swap with new price data)
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Cells.Item(43, 4) "0.8862"
$ws.Cells.Item(43, 5).Value = "  -2.40%  "

# Row 44
Set-TextCell $ws.Cells.Item(44, 4) "0.9998"
$ws.Cells.Item(44, 5).Value = "  -0.02%  "

# Row 45
Set-TextCell $ws.Cells.Item(45, 4) "102.10"
$ws.Cells.Item(45, 5).Value = "  -1.23%  "

# Row 46
Set-TextCell $ws.Cells.Item(46, 4) "2.009.52"
$ws.Cells.Item(46, 5).Value = "  -1.29%  "

# Row 47
Set-TextCell $ws.Cells.Item(47, 4) "0.5179"
$ws.Cells.Item(47, 5).Value = "  -2.80%  "

# Row 48
Set-TextCell $ws.Cells.Item(48, 4) "1.767"
$ws.Cells.Item(48, 5).Value = "  -0.51%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.34%  "

# Row 50
Set-TextCell $ws.Cells.Item(50, 4) "9.313"
$ws.Cells.Item(50, 5).Value = "  +0.18%  "

# Row 51
Set-TextCell $ws.Cells.Item(51, 4) "1.000"
$ws.Cells.Item(51, 5).Value = "  -0.05%  "

